$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two runs of the "In code, of course..." paragraph into a
#    single run (no other formatting change).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "In code, of course, these translate to conditional statements, loops, and using variables.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "In code, of course, these translate to conditional statements, loops, and using variables.",
    2)

# ---------------------------------------------------------------------
# 2) Turn the empty "Output" bookmark paragraph into 3 new bulleted
#    paragraphs (continuing the existing numId=1 list used by the
#    Input/Process sections) followed by two empty, indented
#    paragraphs.
# ---------------------------------------------------------------------

# Find the "Output" heading paragraph; the bookmark paragraph is right
# after it.
$outputIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Output") {
        $outputIdx = $i
        break
    }
}
$idx = $outputIdx + 1

# A paragraph that is already part of the numId=1 list, so the three new
# bullets continue that same list instead of minting a new one.
$listSource = $d.Paragraphs.Item(12)

# --- New bullet paragraph #1 (ilvl 0) -----------------------------------
$bookmarkPara = $d.Paragraphs.Item($idx)
$bookmarkPara.Range.InsertBefore([char]13)
$p1 = $d.Paragraphs.Item($idx)
$p1.Range.Text = "Print out what type of triangle it is (equilateral, isosceles, scalene) "
$p1.Range.ListFormat.List = $listSource.Range.ListFormat.List
$p1.Range.ListFormat.ListLevelNumber = 1
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End - 1
$p1TextOnly = $d.Range($p1Start, $p1End - 1)
$p1TextOnly.Font.Name = "Arial Unicode MS"
$p1SpaceOnly = $d.Range($p1End - 1, $p1End)
$p1SpaceOnly.Font.Name = "Arial Unicode MS"

# --- New bullet paragraph #2 (ilvl 1) -----------------------------------
$bookmarkPara2 = $d.Paragraphs.Item($idx + 1)
$bookmarkPara2.Range.InsertBefore([char]13)
$p2 = $d.Paragraphs.Item($idx + 1)
$p2.Range.Text = "If all angles are the same, two are the same or none of them are, print out the correct triangle"
$p2.Range.ListFormat.List = $listSource.Range.ListFormat.List
$p2.Range.ListFormat.ListLevelNumber = 2

# --- New bullet paragraph #3 (ilvl 2, keeps the _GoBack bookmark) ------
$bookmarkPara3 = $d.Paragraphs.Item($idx + 2)
$insertPos = $bookmarkPara3.Range.Start
$txt3 = "If inputs are not valid, print " + [char]8220 + "error" + [char]8221
$rng3 = $d.Range($insertPos, $insertPos)
$rng3.InsertBefore($txt3)
$p3 = $d.Paragraphs.Item($idx + 2)
$p3.Range.ListFormat.List = $listSource.Range.ListFormat.List
$p3.Range.ListFormat.ListLevelNumber = 3
$p3TextOnly = $d.Range($insertPos, $insertPos + $txt3.Length)
$p3TextOnly.Font.Name = "Arial Unicode MS"

# --- Two new empty, indented paragraphs after the bulleted list --------
$finalEmptyIdx = $idx + 3
$finalEmptyPara = $d.Paragraphs.Item($finalEmptyIdx)
$insertPos2 = $finalEmptyPara.Range.Start
$rng4 = $d.Range($insertPos2, $insertPos2)
$rng4.InsertBefore("ZZTMP1" + [char]13 + "ZZTMP2" + [char]13)

$f1 = $d.Content
$f1.Find.ClearFormatting()
$f1.Find.Execute("ZZTMP1", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
$f2 = $d.Content
$f2.Find.ClearFormatting()
$f2.Find.Execute("ZZTMP2", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$emptyPara1 = $d.Paragraphs.Item($finalEmptyIdx)
$emptyPara1.Range.ParagraphFormat.LeftIndent = 18
$emptyPara2 = $d.Paragraphs.Item($finalEmptyIdx + 1)
$emptyPara2.Range.ParagraphFormat.LeftIndent = 36
